# Applies the dated worksheet update: refreshes the date heading and
# regenerates all 100 arithmetic answers in the 20x5 table.
$d = $word.ActiveDocument

$replacements = @(
    @('2024-09-09 Monday', '2024-09-10 Tuesday'),
    @('94-89=5', '82-36=46'),
    @('70-21=49', '49+38=87'),
    @('8+53=61', '18+38=56'),
    @('35+39=74', '51-45=6'),
    @('90-8=82', '59+18=77'),
    @('28+37=65', '96-7=89'),
    @('87-39=48', '86-19=67'),
    @('38+36=74', '40-28=12'),
    @('49+49=98', '19+38=57'),
    @('73-69=4', '9+73=82'),
    @('82-5=77', '52-35=17'),
    @('6+36=42', '19+23=42'),
    @('75-59=16', '59+14=73'),
    @('95-59=36', '34+9=43'),
    @('61-29=32', '36+25=61'),
    @('81-47=34', '71-7=64'),
    @('36+59=95', '45+37=82'),
    @('74-67=7', '93-88=5'),
    @('40-36=4', '29+46=75'),
    @('4+87=91', '47+17=64'),
    @('90-38=52', '27+36=63'),
    @('79+2=81', '73-26=47'),
    @('23+19=42', '93-56=37'),
    @('46+25=71', '95-38=57'),
    @('87-69=18', '85+9=94'),
    @('69+18=87', '55+16=71'),
    @('8+59=67', '41-25=16'),
    @('76-69=7', '23+48=71'),
    @('92-88=4', '92-35=57'),
    @('72-34=38', '7+45=52'),
    @('70-68=2', '88-69=19'),
    @('64-55=9', '9+52=61'),
    @('96-29=67', '46+39=85'),
    @('26+38=64', '83-28=55'),
    @('90-82=8', '95-79=16'),
    @('84-5=79', '73-15=58'),
    @('70-61=9', '61-52=9'),
    @('71-2=69', '53-7=46'),
    @('84-16=68', '17+17=34'),
    @('94-88=6', '57+39=96'),
    @('16+69=85', '73+18=91'),
    @('80-66=14', '88-59=29'),
    @('83-55=28', '90-83=7'),
    @('28+39=67', '57-28=29'),
    @('69+24=93', '61-53=8'),
    @('36+19=55', '56+8=64'),
    @('71-15=56', '90-73=17'),
    @('35-6=29', '59+14=73'),
    @('42+29=71', '94-28=66'),
    @('19+8=27', '41-5=36'),
    @('68+28=96', '73-64=9'),
    @('87+8=95', '42-35=7'),
    @('84-26=58', '35-27=8'),
    @('84+8=92', '81-43=38'),
    @('94-57=37', '70-19=51'),
    @('18+47=65', '93-24=69'),
    @('65-28=37', '62+29=91'),
    @('9+3=12', '88-69=19'),
    @('80-52=28', '15-7=8'),
    @('62-38=24', '17+47=64'),
    @('80-77=3', '62-38=24'),
    @('43-25=18', '59+7=66'),
    @('26+5=31', '55-38=17'),
    @('33+59=92', '82-36=46'),
    @('81-62=19', '76-68=8'),
    @('90-61=29', '17+76=93'),
    @('43+18=61', '4+68=72'),
    @('90-53=37', '17+65=82'),
    @('50-31=19', '39+3=42'),
    @('53+8=61', '92-79=13'),
    @('29+18=47', '29+46=75'),
    @('18+45=63', '94-69=25'),
    @('81-44=37', '62-54=8'),
    @('94-37=57', '96-47=49'),
    @('64-56=8', '59+3=62'),
    @('45-18=27', '67+26=93'),
    @('32-17=15', '17+34=51'),
    @('36-7=29', '38+54=92'),
    @('3+79=82', '58+19=77'),
    @('39+34=73', '80-41=39'),
    @('49+32=81', '87-68=19'),
    @('9+2=11', '4+49=53'),
    @('83-37=46', '86+6=92'),
    @('85-58=27', '8+66=74'),
    @('5+9=14', '84-19=65'),
    @('52-8=44', '5+68=73'),
    @('8+15=23', '8+63=71'),
    @('60-47=13', '57-49=8'),
    @('49+36=85', '81-13=68'),
    @('84-58=26', '55+38=93'),
    @('78-9=69', '8+39=47'),
    @('64-35=29', '37+35=72'),
    @('55+39=94', '83-26=57'),
    @('20-8=12', '81-77=4'),
    @('91-9=82', '58+19=77'),
    @('68-29=39', '7+24=31'),
    @('54+38=92', '65-58=7'),
    @('31-4=27', '53-29=24'),
    @('85-7=78', '17+75=92'),
    @('38+38=76', '18+17=35'),
)

$notFound = @()
foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        $notFound += $old
    }
}

if ($notFound.Count -gt 0) {
    Write-Host "NOT FOUND: $($notFound -join '; ')"
} else {
    Write-Host "All $($replacements.Count) replacements applied successfully."
}
